# Weekly price-sheet update: a new weekly observation is inserted at the
# top of the data block (row 30, right after the last "fixed" row 29),
# pushing all the subsequent weekly rows down by one. The oldest
# observation that falls off the bottom (old row 58) ends up preserved as
# the new last row (59).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 30; rows 30-58 shift down to 31-59.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with this week's observation.
$ws.Cells.Item(30, 1).Value = 9
$ws.Cells.Item(30, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(30, 3).Value = "Metropolitana"
$ws.Cells.Item(30, 4).Value = 44790
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(30, 6).Value = 100112035
$ws.Cells.Item(30, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 36
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 20000
$ws.Cells.Item(30, 13).Value = 20000
$ws.Cells.Item(30, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(30, 15).Value = "Hijuelas"
$ws.Cells.Item(30, 16).Value = 1333
$ws.Cells.Item(30, 17).Value = 15
$ws.Cells.Item(30, 18).Value = "Hortaliza"
